# Generate Report for Handback
# Regenerate the handback-status report: new handoff/handback GUID file
# names, refreshed correspond-file names and refreshed timestamps.

$wb = $excel.ActiveWorkbook

$newMd1 = "e8b0186b-a8e4-4821-907b-ae5a903d2a7d.md"
$newMd2 = "ffff82d570c0-3e41-4574-a9c4-fbded1a52db1.md"

$newZhCn = "e8b0186b-a8e4-4821-907b-ae5a903d2a7d.de989de1c6f941870db4db0d4dab7f11bb26bad7.zh-cn.xlf"
$newDeDe = "e8b0186b-a8e4-4821-907b-ae5a903d2a7d.de989de1c6f941870db4db0d4dab7f11bb26bad7.de-de.xlf"

$newHandoffZhCn = "2016-03-25 01:27:20"
$newHandbackZhCn = "2016-03-25 01:27:46"
$newHandoffDeDe = "2016-03-25 01:27:25"
$newHandbackDeDe = "2016-03-25 01:27:53"

function Update-Hyperlinks($ws, $map) {
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($map.ContainsKey($addr)) {
            $h.TextToDisplay = $map[$addr]
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $newMd1
$ws1.Range("A3").Value = $newMd2

Update-Hyperlinks $ws1 @{
    '$A$2' = $newMd1
    '$A$3' = $newMd2
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $newMd1
$ws2.Range("D2").Value = $newZhCn
$ws2.Range("E2").Value = $newHandoffZhCn
$ws2.Range("F2").Value = $newMd1
$ws2.Range("G2").Value = $newZhCn
$ws2.Range("H2").Value = $newHandbackZhCn

$ws2.Range("A3").Value = $newMd2
$ws2.Range("D3").Value = $newZhCn
$ws2.Range("E3").Value = $newHandoffZhCn
$ws2.Range("F3").Value = $newMd2
$ws2.Range("G3").Value = $newZhCn
$ws2.Range("H3").Value = $newHandbackZhCn

Update-Hyperlinks $ws2 @{
    '$A$2' = $newMd1
    '$D$2' = $newZhCn
    '$F$2' = $newMd1
    '$G$2' = $newZhCn
    '$A$3' = $newMd2
    '$D$3' = $newZhCn
    '$F$3' = $newMd2
    '$G$3' = $newZhCn
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $newMd1
$ws3.Range("D2").Value = $newDeDe
$ws3.Range("E2").Value = $newHandoffDeDe
$ws3.Range("F2").Value = $newMd1
$ws3.Range("G2").Value = $newDeDe
$ws3.Range("H2").Value = $newHandbackDeDe

$ws3.Range("A3").Value = $newMd2
$ws3.Range("D3").Value = $newDeDe
$ws3.Range("E3").Value = $newHandoffDeDe
$ws3.Range("F3").Value = $newMd2
$ws3.Range("G3").Value = $newDeDe
$ws3.Range("H3").Value = $newHandbackDeDe

Update-Hyperlinks $ws3 @{
    '$A$2' = $newMd1
    '$D$2' = $newDeDe
    '$F$2' = $newMd1
    '$G$2' = $newDeDe
    '$A$3' = $newMd2
    '$D$3' = $newDeDe
    '$F$3' = $newMd2
    '$G$3' = $newDeDe
}
